$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with refreshed crypto data.
# Leading "'" forces text storage so numeric-looking strings (e.g. "214.45")
# are not auto-coerced to numbers by Excel, matching the original inline-string cells.

$ws.Range("D2").Value = "'29.903.81"
$ws.Range("E2").Value = "'  +0.52%  "
$ws.Range("D3").Value = "'1.632.57"
$ws.Range("E3").Value = "'  +1.61%  "
$ws.Range("E4").Value = "'  +0.14%  "
$ws.Range("D5").Value = "'214.45"
$ws.Range("E5").Value = "'  +0.59%  "
$ws.Range("E6").Value = "'  +0.17%  "
$ws.Range("E7").Value = "'  +0.08%  "
$ws.Range("D8").Value = "'28.45"
$ws.Range("E8").Value = "'  +0.62%  "
$ws.Range("E9").Value = "'  +1.70%  "
$ws.Range("E10").Value = "'  +0.76%  "
$ws.Range("D11").Value = "'0.0912"
$ws.Range("E11").Value = "'  +0.24%  "
$ws.Range("D12").Value = "'1.866.96"
$ws.Range("E12").Value = "'  +1.72%  "
$ws.Range("D13").Value = "'1.643.73"
$ws.Range("E13").Value = "'  +2.40%  "
$ws.Range("E14").Value = "'  +2.39%  "
$ws.Range("D15").Value = "'9.27"
$ws.Range("E15").Value = "'  +17.60%  "
$ws.Range("D16").Value = "'29.954.66"
$ws.Range("E16").Value = "'  +0.73%  "
$ws.Range("E17").Value = "'  +2.06%  "
$ws.Range("D18").Value = "'64.04"
$ws.Range("E18").Value = "'  -0.19%  "
$ws.Range("D19").Value = "'242.64"
$ws.Range("E19").Value = "'  +0.14%  "
$ws.Range("E20").Value = "'  +0.30%  "
$ws.Range("E21").Value = "'  +0.05%  "
$ws.Range("D22").Value = "'9.85"
$ws.Range("E22").Value = "'  +4.72%  "
$ws.Range("E23").Value = "'  +2.20%  "
$ws.Range("D24").Value = "'2.14"
$ws.Range("E24").Value = "'  +1.46%  "
$ws.Range("D25").Value = "'157.61"
$ws.Range("E25").Value = "'  +1.60%  "
$ws.Range("D26").Value = "'15.51"
$ws.Range("E26").Value = "'  +0.17%  "
$ws.Range("E27").Value = "'  +1.17%  "
$ws.Range("E28").Value = "'  +2.18%  "
$ws.Range("E29").Value = "'  +0.06%  "
$ws.Range("D30").Value = "'0.0486"
$ws.Range("E30").Value = "'  +1.08%  "
$ws.Range("D31").Value = "'1.11"
$ws.Range("E31").Value = "'  +3.99%  "
$ws.Range("E32").Value = "'  +3.82%  "
$ws.Range("D33").Value = "'3.16"
$ws.Range("E33").Value = "'  -0.87%  "
$ws.Range("D34").Value = "'1.423.02"
$ws.Range("E34").Value = "'  -0.21%  "
$ws.Range("E35").Value = "'  +4.23%  "
$ws.Range("D36").Value = "'1.03"
$ws.Range("E36").Value = "'  +0.11%  "
$ws.Range("D37").Value = "'2.81"
$ws.Range("E37").Value = "'  -4.08%  "
$ws.Range("E38").Value = "'  -0.06%  "
$ws.Range("E39").Value = "'  +0.15%  "
$ws.Range("D40").Value = "'75.68"
$ws.Range("E40").Value = "'  +14.00%  "
$ws.Range("D41").Value = "'0.553"
$ws.Range("E41").Value = "'  +0.95%  "
$ws.Range("E42").Value = "'  +1.81%  "
$ws.Range("D43").Value = "'0.826"
$ws.Range("E43").Value = "'  +1.10%  "
$ws.Range("D44").Value = "'0.0487"
$ws.Range("E44").Value = "'  -1.33%  "
$ws.Range("E45").Value = "'  +0.05%  "
$ws.Range("E46").Value = "'  +3.00%  "
$ws.Range("D47").Value = "'52.86"
$ws.Range("E47").Value = "'  -7.18%  "
$ws.Range("D48").Value = "'1.775.78"
$ws.Range("E49").Value = "'  -0.69%  "
$ws.Range("D50").Value = "'0.0₆0112"
$ws.Range("E50").Value = "'  +8.31%  "
$ws.Range("D51").Value = "'89.70"
$ws.Range("E51").Value = "'  +3.60%  "
